$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 85 <- old row 86 (F:V block)
$ws.Range("F85").Value = "Paderborn"
$ws.Range("G85").Value = 2
$ws.Range("H85").Value = "St. Pauli"
$ws.Range("I85").Value = 2
$ws.Range("J85").Value = 2.81
$ws.Range("K85").Value = "09/10/2023 16:12"
$ws.Range("L85").Value = 3.45
$ws.Range("M85").Value = "21/10/2023 12:45"
$ws.Range("N85").Value = 3.7
$ws.Range("O85").Value = "09/10/2023 16:12"
$ws.Range("P85").Value = 3.84
$ws.Range("Q85").Value = "21/10/2023 12:55"
$ws.Range("R85").Value = 2.42
$ws.Range("S85").Value = "09/10/2023 16:12"
$ws.Range("T85").Value = 2.09
$ws.Range("U85").Value = "21/10/2023 12:52"
$ws.Range("V85").Value = "https://www.betexplorer.com/football/germany/2-bundesliga/paderborn-st-pauli/25tAdyM9/"

# Row 86 <- old row 85 (F:V block)
$ws.Range("F86").Value = "VfL Osnabruck"
$ws.Range("G86").Value = 0
$ws.Range("H86").Value = "Wehen"
$ws.Range("I86").Value = 2
$ws.Range("J86").Value = 2.15
$ws.Range("K86").Value = "10/10/2023 11:42"
$ws.Range("L86").Value = 2.1
$ws.Range("M86").Value = "21/10/2023 12:57"
$ws.Range("N86").Value = 3.81
$ws.Range("O86").Value = "10/10/2023 11:42"
$ws.Range("P86").Value = 3.81
$ws.Range("Q86").Value = "21/10/2023 12:57"
$ws.Range("R86").Value = 3.21
$ws.Range("S86").Value = "10/10/2023 11:42"
$ws.Range("T86").Value = 3.44
$ws.Range("U86").Value = "21/10/2023 12:58"
$ws.Range("V86").Value = "https://www.betexplorer.com/football/germany/2-bundesliga/vfl-osnabruck-wehen/2TOdmXrd/"

# Row 98 <- old row 99 (F:V block)
$ws.Range("F98").Value = "Holstein Kiel"
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = "Nurnberg"
$ws.Range("I98").Value = 2
$ws.Range("J98").Value = 2.01
$ws.Range("K98").Value = "22/10/2023 13:42"
$ws.Range("L98").Value = 2.26
$ws.Range("M98").Value = "29/10/2023 13:14"
$ws.Range("N98").Value = 3.87
$ws.Range("O98").Value = "22/10/2023 13:42"
$ws.Range("P98").Value = 3.71
$ws.Range("Q98").Value = "29/10/2023 13:29"
$ws.Range("R98").Value = 3.62
$ws.Range("S98").Value = "22/10/2023 13:42"
$ws.Range("T98").Value = 3.19
$ws.Range("U98").Value = "29/10/2023 13:29"
$ws.Range("V98").Value = "https://www.betexplorer.com/football/germany/2-bundesliga/holstein-kiel-nurnberg/W2H0ogTE/"

# Row 99 <- old row 100 (F:V block)
$ws.Range("F99").Value = "Magdeburg"
$ws.Range("G99").Value = 1
$ws.Range("H99").Value = "Elversberg"
$ws.Range("I99").Value = 2
$ws.Range("J99").Value = 1.9
$ws.Range("K99").Value = "21/10/2023 14:42"
$ws.Range("L99").Value = 2.08
$ws.Range("M99").Value = "29/10/2023 13:27"
$ws.Range("N99").Value = 4.03
$ws.Range("O99").Value = "21/10/2023 14:42"
$ws.Range("P99").Value = 3.87
$ws.Range("Q99").Value = "29/10/2023 13:27"
$ws.Range("R99").Value = 3.87
$ws.Range("S99").Value = "21/10/2023 14:42"
$ws.Range("T99").Value = 3.46
$ws.Range("U99").Value = "29/10/2023 13:29"
$ws.Range("V99").Value = "https://www.betexplorer.com/football/germany/2-bundesliga/magdeburg-elversberg/SUA9qXcR/"

# Row 100 <- old row 98 (F:V block)
$ws.Range("F100").Value = "Wehen"
$ws.Range("G100").Value = 1
$ws.Range("H100").Value = "Hansa Rostock"
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 2.22
$ws.Range("K100").Value = "22/10/2023 13:42"
$ws.Range("L100").Value = 2.48
$ws.Range("M100").Value = "29/10/2023 13:26"
$ws.Range("N100").Value = 3.65
$ws.Range("O100").Value = "22/10/2023 13:42"
$ws.Range("P100").Value = 3.42
$ws.Range("Q100").Value = "29/10/2023 13:26"
$ws.Range("R100").Value = 3.28
$ws.Range("S100").Value = "22/10/2023 13:42"
$ws.Range("T100").Value = 3.02
$ws.Range("U100").Value = "29/10/2023 13:25"
$ws.Range("V100").Value = "https://www.betexplorer.com/football/germany/2-bundesliga/wehen-hansa-rostock/x40qwVze/"

# Row 101 <- old row 102 (F:V block)
$ws.Range("F101").Value = "Elversberg"
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = "St. Pauli"
$ws.Range("I101").Value = 2
$ws.Range("J101").Value = 3.18
$ws.Range("K101").Value = "29/10/2023 13:42"
$ws.Range("L101").Value = 3.25
$ws.Range("M101").Value = "03/11/2023 18:29"
$ws.Range("N101").Value = 3.79
$ws.Range("O101").Value = "29/10/2023 13:42"
$ws.Range("P101").Value = 3.67
$ws.Range("Q101").Value = "03/11/2023 18:29"
$ws.Range("R101").Value = 2.17
$ws.Range("S101").Value = "29/10/2023 13:42"
$ws.Range("T101").Value = 2.23
$ws.Range("U101").Value = "03/11/2023 18:29"
$ws.Range("V101").Value = "https://www.betexplorer.com/football/germany/2-bundesliga/elversberg-st-pauli/dEGIiUkl/"

# Row 102 <- old row 101 (F:V block)
$ws.Range("F102").Value = "Dusseldorf"
$ws.Range("G102").Value = 1
$ws.Range("H102").Value = "Wehen"
$ws.Range("I102").Value = 3
$ws.Range("J102").Value = 1.56
$ws.Range("K102").Value = "29/10/2023 13:42"
$ws.Range("L102").Value = 1.55
$ws.Range("M102").Value = "03/11/2023 18:29"
$ws.Range("N102").Value = 4.69
$ws.Range("O102").Value = "29/10/2023 13:42"
$ws.Range("P102").Value = 4.66
$ws.Range("Q102").Value = "03/11/2023 18:29"
$ws.Range("R102").Value = 5.34
$ws.Range("S102").Value = "29/10/2023 13:42"
$ws.Range("T102").Value = 5.83
$ws.Range("U102").Value = "03/11/2023 18:29"
$ws.Range("V102").Value = "https://www.betexplorer.com/football/germany/2-bundesliga/dusseldorf-wehen/bJlhy957/"

# Row 103 <- old row 105 (F:V block)
$ws.Range("F103").Value = "Kaiserslautern"
$ws.Range("G103").Value = 0
$ws.Range("H103").Value = "Greuther Furth"
$ws.Range("I103").Value = 2
$ws.Range("J103").Value = 2.39
$ws.Range("K103").Value = "28/10/2023 20:43"
$ws.Range("L103").Value = 2.32
$ws.Range("M103").Value = "04/11/2023 12:59"
$ws.Range("N103").Value = 3.57
$ws.Range("O103").Value = "28/10/2023 20:43"
$ws.Range("P103").Value = 3.58
$ws.Range("Q103").Value = "04/11/2023 12:58"
$ws.Range("R103").Value = 3
$ws.Range("S103").Value = "28/10/2023 20:43"
$ws.Range("T103").Value = 3.16
$ws.Range("U103").Value = "04/11/2023 12:59"
$ws.Range("V103").Value = "https://www.betexplorer.com/football/germany/2-bundesliga/kaiserslautern-greuther-furth/4GyvaCcE/"

# Row 105 <- old row 103 (F:V block)
$ws.Range("F105").Value = "VfL Osnabruck"
$ws.Range("G105").Value = 1
$ws.Range("H105").Value = "Holstein Kiel"
$ws.Range("I105").Value = 1
$ws.Range("J105").Value = 2.43
$ws.Range("K105").Value = "29/10/2023 13:42"
$ws.Range("L105").Value = 2.83
$ws.Range("M105").Value = "04/11/2023 12:59"
$ws.Range("N105").Value = 3.85
$ws.Range("O105").Value = "29/10/2023 13:42"
$ws.Range("P105").Value = 3.66
$ws.Range("Q105").Value = "04/11/2023 12:54"
$ws.Range("R105").Value = 2.72
$ws.Range("S105").Value = "29/10/2023 13:42"
$ws.Range("T105").Value = 2.5
$ws.Range("U105").Value = "04/11/2023 12:54"
$ws.Range("V105").Value = "https://www.betexplorer.com/football/germany/2-bundesliga/vfl-osnabruck-holstein-kiel/EkBNjl5f/"

# Row 113 <- old row 114 (F:V block)
$ws.Range("F113").Value = "Paderborn"
$ws.Range("G113").Value = 1
$ws.Range("H113").Value = "Nurnberg"
$ws.Range("I113").Value = 3
$ws.Range("J113").Value = 1.95
$ws.Range("K113").Value = "05/11/2023 13:42"
$ws.Range("L113").Value = 2.27
$ws.Range("M113").Value = "11/11/2023 12:59"
$ws.Range("N113").Value = 3.96
$ws.Range("O113").Value = "05/11/2023 13:42"
$ws.Range("P113").Value = 3.86
$ws.Range("Q113").Value = "11/11/2023 12:58"
$ws.Range("R113").Value = 3.73
$ws.Range("S113").Value = "05/11/2023 13:42"
$ws.Range("T113").Value = 3.04
$ws.Range("U113").Value = "11/11/2023 12:59"
$ws.Range("V113").Value = "https://www.betexplorer.com/football/germany/2-bundesliga/paderborn-nurnberg/0riICVRD/"

# Row 114 <- old row 113 (F:V block)
$ws.Range("F114").Value = "Braunschweig"
$ws.Range("G114").Value = 3
$ws.Range("H114").Value = "VfL Osnabruck"
$ws.Range("I114").Value = 2
$ws.Range("J114").Value = 2.16
$ws.Range("K114").Value = "05/11/2023 13:42"
$ws.Range("L114").Value = 2.45
$ws.Range("M114").Value = "11/11/2023 12:59"
$ws.Range("N114").Value = 3.79
$ws.Range("O114").Value = "05/11/2023 13:42"
$ws.Range("P114").Value = 3.38
$ws.Range("Q114").Value = "11/11/2023 12:59"
$ws.Range("R114").Value = 3.19
$ws.Range("S114").Value = "05/11/2023 13:42"
$ws.Range("T114").Value = 3.09
$ws.Range("U114").Value = "11/11/2023 12:59"
$ws.Range("V114").Value = "https://www.betexplorer.com/football/germany/2-bundesliga/braunschweig-vfl-osnabruck/htWd5mKm/"
# New row 115
$ws.Range("A114").Copy()
$ws.Range("A115").PasteSpecial(-4122)
$ws.Range("E114").Copy()
$ws.Range("E115").PasteSpecial(-4122)
$ws.Range("A115").Value = 114
$ws.Range("B115").Value = "germany"
$ws.Range("C115").Value = "2-bundesliga"
$ws.Range("D115").Value = "2023-2024"
$ws.Range("E115").Value = 45241.85416666666
$ws.Range("F115").Value = "Hertha Berlin"
$ws.Range("G115").Value = 2
$ws.Range("H115").Value = "Karlsruher SC"
$ws.Range("I115").Value = 2
$ws.Range("J115").Value = 2.06
$ws.Range("K115").Value = "05/11/2023 13:42"
$ws.Range("L115").Value = 1.91
$ws.Range("M115").Value = "11/11/2023 20:28"
$ws.Range("N115").Value = 3.93
$ws.Range("O115").Value = "05/11/2023 13:42"
$ws.Range("P115").Value = 4.08
$ws.Range("Q115").Value = "11/11/2023 20:28"
$ws.Range("R115").Value = 3.41
$ws.Range("S115").Value = "05/11/2023 13:42"
$ws.Range("T115").Value = 3.83
$ws.Range("U115").Value = "11/11/2023 20:28"
$ws.Range("V115").Value = "https://www.betexplorer.com/football/germany/2-bundesliga/hertha-berlin-karlsruher/n39VlSZ6/"
